$wb = $excel.ActiveWorkbook

# --- Sheet1: rearrange scene-quest ids 42040001-8 -> 42020001-8 (rows 41-48) ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Cells.Item(41, 1).Value = 42020001
$ws1.Cells.Item(42, 1).Value = 42020002
$ws1.Cells.Item(43, 1).Value = 42020003
$ws1.Cells.Item(44, 1).Value = 42020004
$ws1.Cells.Item(45, 1).Value = 42020005
$ws1.Cells.Item(46, 1).Value = 42020006
$ws1.Cells.Item(47, 1).Value = 42020007
$ws1.Cells.Item(48, 1).Value = 42020008

# --- MapSet: rearrange scene-quest ids (rows 4-33) ---
$ws2 = $wb.Worksheets.Item("MapSet")

$ws2.Cells.Item(4, 1).Value = 42110001
$ws2.Cells.Item(5, 1).Value = 42110002
$ws2.Cells.Item(6, 1).Value = 42110003
$ws2.Cells.Item(7, 1).Value = 42110004
$ws2.Cells.Item(8, 1).Value = 42110005

$ws2.Cells.Item(9, 1).Value = 42120001
$ws2.Cells.Item(10, 1).Value = 42120002
$ws2.Cells.Item(11, 1).Value = 42120003
$ws2.Cells.Item(12, 1).Value = 42120004
$ws2.Cells.Item(13, 1).Value = 42120005
$ws2.Cells.Item(14, 1).Value = 42120006
$ws2.Cells.Item(15, 1).Value = 42120007
$ws2.Cells.Item(16, 1).Value = 42120008
$ws2.Cells.Item(17, 1).Value = 42120009
$ws2.Cells.Item(18, 1).Value = 42120010
$ws2.Cells.Item(19, 1).Value = 42120011
$ws2.Cells.Item(20, 1).Value = 42120012
$ws2.Cells.Item(21, 1).Value = 42120013
$ws2.Cells.Item(22, 1).Value = 42120014
$ws2.Cells.Item(23, 1).Value = 42120015
$ws2.Cells.Item(24, 1).Value = 42120016
$ws2.Cells.Item(25, 1).Value = 42120017
$ws2.Cells.Item(26, 1).Value = 42120018
$ws2.Cells.Item(27, 1).Value = 42120019

$ws2.Cells.Item(28, 1).Value = 42130001
$ws2.Cells.Item(29, 1).Value = 42130002
$ws2.Cells.Item(30, 1).Value = 42130003
$ws2.Cells.Item(31, 1).Value = 42130004
$ws2.Cells.Item(32, 1).Value = 42130005
$ws2.Cells.Item(33, 1).Value = 42130006

# --- Update the saved view/selection state on both sheets ---
# MapSet: selection moves from J1:J3 to B29 (pane/topLeftCell stays at A4)
$ws2.Activate()
$ws2.Range("B29").Select()

# Sheet1: frozen-pane scroll moves from A16 to A4, selection moves from E38 to B17
# Re-activating Sheet1 last keeps it the visible/active tab, matching the source file.
$ws1.Activate()
$ws1.Range("B17").Select()
